$d = $word.ActiveDocument

# Locate the "Design Skiplist" hyperlink text and collapse the found
# range to its end so we can append a new run right after it (still
# inside the same paragraph, after the closing </w:hyperlink>).
$rng = $d.Content
$found = $rng.Find.Execute("Design Skiplist", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found -eq $false) {
    throw "Could not find 'Design Skiplist' in the document"
}

$rng.Collapse(0)  # wdCollapseEnd

# Match the Arial formatting used elsewhere for trailing annotation
# text (e.g. " (solved)") appended after a hyperlink in this list.
$rng.Font.Name = "Arial"
$rng.Font.NameAscii = "Arial"
$rng.Font.NameFarEast = "Arial"
$rng.Font.Color = 0          # auto
$rng.Font.Spacing = 0
$rng.Font.Position = 0
$rng.Font.Size = 11
$rng.Font.Underline = 0

$rng.InsertAfter(" - hard1")
